# 自动更新Excel文件 - 2025-11-12 23:12:50
#
# Daily "remaining days" countdown rollover:
#   - Column D = total cycle length in days ("总天")
#   - Column E = days remaining ("剩余")
#   - Column F = cycle start date, yyyymmdd integer ("开始时间")
#
# Each day E decrements by 1. When a row's remaining count would drop to 0
# (i.e. it was 1 before today's update), the cycle renews: the start date
# (F) resets to "today" and the remaining count (E) resets to the full
# cycle length (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251113

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $eVal = $eCell.Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    if ($eVal -eq $null) {
        continue
    }

    # Skip rows whose start date isn't a clean 8-digit yyyymmdd stamp
    # (data-entry glitch) - those were left untouched by the refresh.
    $fText = "$fVal"
    if ($fText.Length -ne 8) {
        continue
    }

    if ($eVal -eq 1) {
        # Remaining days would hit zero today -> renew the cycle.
        $dVal = $ws.Cells.Item($r, 4).Value2
        $ws.Cells.Item($r, 6).Value2 = $today
        $eCell.Value2 = $dVal
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
